# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before the
#    existing "2022-Q3" sheet), populated with the Q4 fund-holding table.
# 2. Update the "总计" (summary) sheet: the new Q4 row is prepended and all
#    the other quarters shift down one row, with a brand-new last row for
#    the quarter that used to be in row 4 (2022-Q1).
# 3. Restore the original active-sheet/tab selection (adding a sheet makes
#    it the active one in Excel, but the workbook was last viewed on
#    "2022-Q1").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q4" sheet, inserted before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q4.Cells.Item(1, $i + 2)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# index(num), 基金代码(text), 基金名称(text), 基金规模(text), 股票总仓位(text),
# 仓位占比(text), 持有市值(text), 仓位排名(num)
$q4rows = @(
    @(0, "000522", "华润元大信息传媒科技混合", "1.38", "62.01", "6.62", "0.0914", 1),
    @(1, "012075", "易方达稳健添利混合A",     "1.39", "47.10", "3.55", "0.0493", 6),
    @(2, "012076", "易方达稳健添利混合C",     "0.72", "47.10", "3.55", "0.0256", 6)
)

foreach ($r in $q4rows) {
    $row = [int]$r[0] + 2
    $a = $q4.Cells.Item($row, 1)
    $a.Value = $r[0]
    $a.Font.Bold = $true
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4160
    $a.Borders.LineStyle = 1

    for ($col = 2; $col -le 7; $col++) {
        $textCell = $q4.Cells.Item($row, $col)
        $textCell.NumberFormat = "@"
        $textCell.Value = $r[$col - 1]
    }

    $q4.Cells.Item($row, 8).Value = $r[7]
}

# ---------------------------------------------------------------------
# 2) Update "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.17

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 9
$total.Range("D3").Value = 0.51

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.15

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 5
$total.Range("D5").Value = 0.18

$total.Range("A5").Font.Bold = $true
$total.Range("A5").HorizontalAlignment = -4108
$total.Range("A5").VerticalAlignment = -4160
$total.Range("A5").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3) Restore the active sheet/tab to "2022-Q1" (last sheet), since adding
#    the new worksheet above shifted Excel's active-tab selection.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Activate()
$q1.Select()
